# Applies cached-value corrections to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled market-data refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 255.44827
$ws.Range("I33").Value = 212.8
$ws.Range("K33").Value = 212.8
$ws.Range("M33").Value = 16.19999999999999
$ws.Range("H40").Value = 941.8393
$ws.Range("I40").Value = 922.9474
$ws.Range("J40").Value = 981.7222
$ws.Range("K40").Value = 922.9474
$ws.Range("L40").Value = 981.7222
$ws.Range("M40").Value = -747.9474
$ws.Range("N40").Value = -1331.7222
$ws.Range("H100").Value = 1646.1538
$ws.Range("I100").Value = 1357.1428
$ws.Range("J100").Value = 1983.3334
$ws.Range("K100").Value = 1357.1428
$ws.Range("L100").Value = 1983.3334
$ws.Range("M100").Value = -816.1428000000001
$ws.Range("N100").Value = -3065.3334
$ws.Range("H132").Value = 440188.8
$ws.Range("I132").Value = 561907.9399999999
$ws.Range("K132").Value = 1685723.82
$ws.Range("M132").Value = -1683193.82
$ws.Range("H138").Value = 2319.6416
$ws.Range("I138").Value = 1900.2916
$ws.Range("J138").Value = 2666.6897
$ws.Range("K138").Value = 5700.8748
$ws.Range("L138").Value = 8000.0691
$ws.Range("M138").Value = -560.8747999999996
$ws.Range("N138").Value = -18280.0691
$ws.Range("H139").Value = 54500
$ws.Range("J139").Value = 54500
$ws.Range("L139").Value = 54500
$ws.Range("N139").Value = -64780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 10000000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H12").Value = 3366.6667
$ws.Range("I12").Value = 600
$ws.Range("J12").Value = 4750
$ws.Range("K12").Value = 600
$ws.Range("L12").Value = 4750
$ws.Range("M12").Value = -427
$ws.Range("N12").Value = -5096
$ws.Range("H32").Value = 4743.365
$ws.Range("I32").Value = 5368.826
$ws.Range("K32").Value = 5368.826
$ws.Range("M32").Value = -5081.826
$ws.Range("H61").Value = 4480.8887
$ws.Range("I61").Value = 2600
$ws.Range("J61").Value = 6832
$ws.Range("K61").Value = 2600
$ws.Range("L61").Value = 6832
$ws.Range("M61").Value = -2388
$ws.Range("N61").Value = -7256
$ws.Range("H74").Value = 3534.5334
$ws.Range("I74").Value = 1040.7142
$ws.Range("J74").Value = 4660.7744
$ws.Range("K74").Value = 1040.7142
$ws.Range("L74").Value = 4660.7744
$ws.Range("M74").Value = -166.7141999999999
$ws.Range("N74").Value = -6408.7744
$ws.Range("H77").Value = 3534.5334
$ws.Range("I77").Value = 1040.7142
$ws.Range("J77").Value = 4660.7744
$ws.Range("K77").Value = 5203.571
$ws.Range("L77").Value = 23303.872
$ws.Range("M77").Value = -835.5709999999999
$ws.Range("N77").Value = -32039.872
$ws.Range("H82").Value = 39900
$ws.Range("J82").Value = 39900
$ws.Range("L82").Value = 39900
$ws.Range("N82").Value = -40622
$ws.Range("H85").Value = 39900
$ws.Range("J85").Value = 39900
$ws.Range("L85").Value = 39900
$ws.Range("N85").Value = -42396
$ws.Range("H86").Value = 20649.5
$ws.Range("J86").Value = 39999
$ws.Range("L86").Value = 39999
$ws.Range("N86").Value = -42371
$ws.Range("H89").Value = 20649.5
$ws.Range("J89").Value = 39999
$ws.Range("L89").Value = 119997
$ws.Range("N89").Value = -131853
$ws.Range("H96").Value = 34344
$ws.Range("J96").Value = 34344
$ws.Range("L96").Value = 34344
$ws.Range("N96").Value = -39836
$ws.Range("H132").Value = 3536.6365
$ws.Range("I132").Value = 3341.647
$ws.Range("K132").Value = 10024.941
$ws.Range("M132").Value = -7494.940999999999
$ws.Range("H136").Value = 4480.8887
$ws.Range("I136").Value = 2600
$ws.Range("J136").Value = 6832
$ws.Range("K136").Value = 7800
$ws.Range("L136").Value = 20496
$ws.Range("M136").Value = -5250
$ws.Range("N136").Value = -25596

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 58739.332
$ws.Range("I134").Value = 64764.527
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 194293.581
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -191758.581
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 13000
$ws.Range("J28").Value = 13000
$ws.Range("L28").Value = 13000
$ws.Range("N28").Value = -13490
$ws.Range("H31").Value = 1638
$ws.Range("I31").Value = 906.82355
$ws.Range("J31").Value = 2466.6667
$ws.Range("K31").Value = 906.82355
$ws.Range("L31").Value = 2466.6667
$ws.Range("M31").Value = -611.82355
$ws.Range("N31").Value = -3056.6667
$ws.Range("H34").Value = 1638
$ws.Range("I34").Value = 906.82355
$ws.Range("J34").Value = 2466.6667
$ws.Range("K34").Value = 906.82355
$ws.Range("L34").Value = 2466.6667
$ws.Range("M34").Value = -704.82355
$ws.Range("N34").Value = -2870.6667
$ws.Range("H58").Value = 2882.4
$ws.Range("I58").Value = 2936.2856
$ws.Range("J58").Value = 2853.3845
$ws.Range("K58").Value = 2936.2856
$ws.Range("L58").Value = 2853.3845
$ws.Range("M58").Value = -2733.2856
$ws.Range("N58").Value = -3259.3845
$ws.Range("H125").Value = 50393.2
$ws.Range("J125").Value = 50393.2
$ws.Range("L125").Value = 50393.2
$ws.Range("N125").Value = -55313.2
$ws.Range("H136").Value = 2882.4
$ws.Range("I136").Value = 2936.2856
$ws.Range("J136").Value = 2853.3845
$ws.Range("K136").Value = 8808.856800000001
$ws.Range("L136").Value = 8560.1535
$ws.Range("M136").Value = -6258.856800000001
$ws.Range("N136").Value = -13660.1535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 45755.09
$ws.Range("I107").Value = 50154.1
$ws.Range("J107").Value = 42089.25
$ws.Range("K107").Value = 150462.3
$ws.Range("L107").Value = 126267.75
$ws.Range("M107").Value = -148542.3
$ws.Range("N107").Value = -130107.75
$ws.Range("H117").Value = 1689.6
$ws.Range("I117").Value = 792.8570999999999
$ws.Range("J117").Value = 2474.25
$ws.Range("K117").Value = 2378.5713
$ws.Range("L117").Value = 7422.75
$ws.Range("M117").Value = 1063.4287
$ws.Range("N117").Value = -14306.75
$ws.Range("H121").Value = 33334432
$ws.Range("I121").Value = 50
$ws.Range("J121").Value = 38462800
$ws.Range("K121").Value = 150
$ws.Range("L121").Value = 115388400
$ws.Range("M121").Value = 1160
$ws.Range("N121").Value = -115391020
$ws.Range("H129").Value = 1137.9
$ws.Range("J129").Value = 1500
$ws.Range("L129").Value = 4500
$ws.Range("N129").Value = -14500
$ws.Range("H132").Value = 45455384
$ws.Range("I132").Value = 62500584
$ws.Range("K132").Value = 562505256
$ws.Range("M132").Value = -562502726

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 600
$ws.Range("J13").Value = 600
$ws.Range("L13").Value = 600
$ws.Range("N13").Value = -878

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 368.36365
$ws.Range("I22").Value = 335.7143
$ws.Range("J22").Value = 425.5
$ws.Range("K22").Value = 335.7143
$ws.Range("L22").Value = 425.5
$ws.Range("M22").Value = -40.71429999999998
$ws.Range("N22").Value = -1015.5
$ws.Range("H27").Value = 368.36365
$ws.Range("I27").Value = 335.7143
$ws.Range("J27").Value = 425.5
$ws.Range("K27").Value = 335.7143
$ws.Range("L27").Value = 425.5
$ws.Range("M27").Value = -228.7143
$ws.Range("N27").Value = -639.5
$ws.Range("H62").Value = 24999.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 24999.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 24999.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -26247.5
$ws.Range("H65").Value = 24999.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 24999.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 74998.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -81238.5
$ws.Range("H94").Value = 24900
$ws.Range("J94").Value = 24900
$ws.Range("L94").Value = 24900
$ws.Range("N94").Value = -26252
$ws.Range("H122").Value = 7210.909
$ws.Range("I122").Value = 8168.8887
$ws.Range("J122").Value = 2900
$ws.Range("K122").Value = 24506.6661
$ws.Range("L122").Value = 8700
$ws.Range("M122").Value = -22056.6661
$ws.Range("N122").Value = -13600
$ws.Range("H132").Value = 189183.67
$ws.Range("I132").Value = 280776
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 842328
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -839798
$ws.Range("N132").Value = -23057
$ws.Range("H139").Value = 54465
$ws.Range("J139").Value = 54465
$ws.Range("L139").Value = 54465
$ws.Range("N139").Value = -64745

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 12501500
$ws.Range("I7").Value = 25000250
$ws.Range("J7").Value = 2750
$ws.Range("K7").Value = 25000250
$ws.Range("L7").Value = 2750
$ws.Range("M7").Value = -25000137
$ws.Range("N7").Value = -2976
